$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "Transmission Map Generator" table values (rows 8-9) ---
$ws.Range("G8").Value = 33.608440000000002
$ws.Range("H8").Value = 0.00044999999999999999
$ws.Range("I8").Value = 0.97160000000000002

$ws.Range("G9").Value = 29.921849999999999
$ws.Range("H9").Value = 0.00103999999999999991
$ws.Range("I9").Value = 0.95730000000000004

# --- Add new model version rows (10-12) to the "Transmission Map Generator" table ---
$ws.Range("F10").Value = "V1.04.5"
$ws.Range("G10").Value = 32.728729999999999
$ws.Range("H10").Value = 0.00055000000000000003
$ws.Range("I10").Value = 0.96704999999999997

$ws.Range("F11").Value = "V1.04.6"
$ws.Range("G11").Value = 29.98086
$ws.Range("H11").Value = 0.00102000000000000007
$ws.Range("I11").Value = 0.93915999999999999

$ws.Range("F12").Value = "V1.04.7"
$ws.Range("G12").Value = 29.826049999999999
$ws.Range("H12").Value = 0.00105999999999999996
$ws.Range("I12").Value = 0.93645

# --- Add new row for "RESIDE-OTS Performance" table (columns W/X) ---
$ws.Range("W4").Value = 22.3033
$ws.Range("W4").NumberFormat = "0.00000"
$ws.Range("X4").Value = 0.96426000000000001

# --- Update selected cell to reflect where the author last left off ---
[void]$ws.Range("J12").Select()

"done"
